# "load start samosbor fix" — update the task list:
#  - Row 6 ("add single-recipe crafting system") gets a longer description
#    and its status flips from "Не сделано" (not done) to "В процессе"
#    (in progress) with a "0.6.0" target-version tag.
#  - Row 9 ("Добавить чат (say)") and Row 15 ("Сделать усложнение врагов...")
#    also flip from "Не сделано" to "В процессе" with a "0.6.0" tag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: rewrite the task text and mark it "in progress" / "0.6.0" ---
$ws.Range("A6").Value = "добавить систему рецептов из одного предмета и просмотр крафта определенного предмета"

$ws.Range("D4").Copy()
$ws.Range("B6").PasteSpecial(-4122)          # xlPasteFormats: copy "В процессе" style
$ws.Range("B6").Value = "В процессе"

$ws.Range("C1").Copy()
$ws.Range("C6").PasteSpecial(-4122)          # xlPasteFormats: copy the "0.6.0" tag style
$ws.Range("C6").Value = "0.6.0"

$ws.Rows.Item(6).RowHeight = 45

# --- Row 9: same status flip ---
$ws.Range("D4").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "В процессе"

$ws.Range("C1").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = "0.6.0"

# --- Row 15: same status flip ---
$ws.Range("D4").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "В процессе"

$ws.Range("C1").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "0.6.0"

$excel.CutCopyMode = $false
